$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$numericDCells = @("D8","D9","D16","D18","D19","D22","D24","D26","D29","D31","D33","D34","D41","D44","D46","D47","D50","D51")
foreach ($ref in $numericDCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = "29.603.85"
$ws.Range("E2").Value = "  +3.30%  "
$ws.Range("D3").Value = "1.606.17"
$ws.Range("E3").Value = "  +2.71%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("E5").Value = "  +0.94%  "
$ws.Range("E6").Value = "  +2.67%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").Value = "26.81"
$ws.Range("E8").Value = "  +7.55%  "
$ws.Range("D9").Value = "43.58"
$ws.Range("E9").Value = "  -1.28%  "
$ws.Range("E10").Value = "  +2.21%  "
$ws.Range("E11").Value = "  +2.58%  "
$ws.Range("E12").Value = "  +1.29%  "
$ws.Range("D13").Value = "1.834.28"
$ws.Range("E13").Value = "  +2.56%  "
$ws.Range("D14").Value = "1.594.43"
$ws.Range("E14").Value = "  +2.00%  "
$ws.Range("D15").Value = "29.604.51"
$ws.Range("E15").Value = "  +3.14%  "
$ws.Range("D16").Value = "0.537"
$ws.Range("E16").Value = "  +3.86%  "
$ws.Range("D18").Value = "63.43"
$ws.Range("E18").Value = "  +3.16%  "
$ws.Range("D19").Value = "240.53"
$ws.Range("E19").Value = "  +5.57%  "
$ws.Range("E20").Value = "  +3.76%  "
$ws.Range("D21").Value = "0.0₃0693"
$ws.Range("E21").Value = "  +1.64%  "
$ws.Range("D22").Value = "0.999"
$ws.Range("E23").Value = "  +1.58%  "
$ws.Range("D24").Value = "9.21"
$ws.Range("E24").Value = "  +2.01%  "
$ws.Range("E25").Value = "  +0.51%  "
$ws.Range("D26").Value = "154.40"
$ws.Range("E26").Value = "  +1.49%  "
$ws.Range("E27").Value = "  +2.68%  "
$ws.Range("E28").Value = "  +3.39%  "
$ws.Range("D29").Value = "6.40"
$ws.Range("E29").Value = "  +2.70%  "
$ws.Range("E30").Value = "  -0.06%  "
$ws.Range("D31").Value = "0.0473"
$ws.Range("E31").Value = "  +3.27%  "
$ws.Range("E32").Value = "  +0.76%  "
$ws.Range("D33").Value = "3.23"
$ws.Range("E33").Value = "  +1.63%  "
$ws.Range("D34").Value = "3.11"
$ws.Range("E34").Value = "  +3.98%  "
$ws.Range("D35").Value = "1.407.86"
$ws.Range("E35").Value = "  +0.55%  "
$ws.Range("E36").Value = "  +0.89%  "
$ws.Range("E37").Value = "  +5.30%  "
$ws.Range("E38").Value = "  +5.77%  "
$ws.Range("E39").Value = "  +0.11%  "
$ws.Range("E40").Value = "  +2.39%  "
$ws.Range("D41").Value = "0.538"
$ws.Range("E41").Value = "  +4.02%  "
$ws.Range("E42").Value = "  +2.33%  "
$ws.Range("E43").Value = "  +7.69%  "
$ws.Range("D44").Value = "53.83"
$ws.Range("E44").Value = "  +27.01%  "
$ws.Range("E45").Value = "  +3.76%  "
$ws.Range("D46").Value = "0.999"
$ws.Range("E46").Value = "  -0.09%  "
$ws.Range("D47").Value = "65.92"
$ws.Range("E47").Value = "  +3.17%  "
$ws.Range("E48").Value = "  +1.20%  "
$ws.Range("D49").Value = "1.745.80"
$ws.Range("E49").Value = "  +2.53%  "
$ws.Range("D50").Value = "0.863"
$ws.Range("E50").Value = "  -0.20%  "
$ws.Range("D51").Value = "86.61"
$ws.Range("E51").Value = "  +2.18%  "

foreach ($ref in $numericDCells) {
    $ws.Range($ref).Style = "Normal"
}
